$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21 ("10" / JC row): JE -> JR, JNE -> JA, JB removed
$ws.Range("C21").Value = "JR"
$ws.Range("D21").Value = "JA"
$ws.Range("E21").ClearContents()

# Row 22 ("11" row): JR/JA removed (merged up into row21), CR and EXT remain
$ws.Range("B22").ClearContents()
$ws.Range("C22").ClearContents()

# Update the active selection to C21
$ws.Range("C21").Select()
